$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (G1/H1) - set values then copy the bold/border header style from F1
$ws.Range("G1").Value = "Temp(c)"
$ws.Range("H1").Value = "Temp(f)"
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# Row 2 - existing Pune City West row, A2 becomes a plain numeric pincode
$ws.Range("A2").Value = 411007
$ws.Range("B2").Value = "Pune City West"

# Row 3 - Udupi
$ws.Range("A3").Value = 576104
$ws.Range("B3").Value = "Udupi"
$ws.Range("E3").Value = 13.35
$ws.Range("F3").Value = 74.75
$ws.Range("G3").Value = 25
$ws.Range("H3").Value = 77

# Row 4 - Gurgaon
$ws.Range("A4").Value = 122101
$ws.Range("B4").Value = "Gurgaon"
$ws.Range("E4").Value = 28.47
$ws.Range("F4").Value = 77.03
$ws.Range("G4").Value = 14
$ws.Range("H4").Value = 57.2

# Row 5 - Gurgaon
$ws.Range("A5").Value = 122003
$ws.Range("B5").Value = "Gurgaon"
$ws.Range("E5").Value = 28.47
$ws.Range("F5").Value = 77.03
$ws.Range("G5").Value = 14
$ws.Range("H5").Value = 57.2

# Row 6 - Pune City West again, but this time the pincode is stored as text
$ws.Range("A6").Value = "'411007"
$ws.Range("A6").Style = $ws.Range("A2").Style
$ws.Range("B6").Value = "Pune City West"
$ws.Range("E6").Value = 18.53
$ws.Range("F6").Value = 73.87
$ws.Range("G6").Value = 17.4
$ws.Range("H6").Value = 63.4
